# Weekly Fruta/Hortaliza update: insert the latest week's two new price
# records (Primera / Segunda) for Maracuyá at the top of the data block
# (rows 34-35), pushing all the older rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 34 and 35 (everything below shifts down by 2).
$ws.Range("A34:A35").EntireRow.Insert()

# Common / constant column values shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100108
$producto  = "Tropicales y subtropicales"
$categoriaId = 100108003
$categoria = "Maracuyá"
$variedad  = "Sin especificar"
$unidad    = "$/caja 20 kilos"
$origen    = "Región de Arica y Parinacota"
$kgUnidad  = 20

# --- New row 34: Primera ---
$ws.Cells.Item(34, 1).Value  = $mercadoId
$ws.Cells.Item(34, 2).Value  = $mercado
$ws.Cells.Item(34, 3).Value  = $region
$ws.Cells.Item(34, 4).Value  = 44648
$ws.Cells.Item(34, 5).Value  = $codreg
$ws.Cells.Item(34, 6).Value  = $tipo
$ws.Cells.Item(34, 7).Value  = $productoId
$ws.Cells.Item(34, 8).Value  = $producto
$ws.Cells.Item(34, 9).Value  = $categoriaId
$ws.Cells.Item(34, 10).Value = $categoria
$ws.Cells.Item(34, 11).Value = $variedad
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 120
$ws.Cells.Item(34, 14).Value = 17000
$ws.Cells.Item(34, 15).Value = 18000
$ws.Cells.Item(34, 16).Value = 17500
$ws.Cells.Item(34, 17).Value = $unidad
$ws.Cells.Item(34, 18).Value = $origen
$ws.Cells.Item(34, 19).Value = 875
$ws.Cells.Item(34, 20).Value = $kgUnidad

# --- New row 35: Segunda ---
$ws.Cells.Item(35, 1).Value  = $mercadoId
$ws.Cells.Item(35, 2).Value  = $mercado
$ws.Cells.Item(35, 3).Value  = $region
$ws.Cells.Item(35, 4).Value  = 44648
$ws.Cells.Item(35, 5).Value  = $codreg
$ws.Cells.Item(35, 6).Value  = $tipo
$ws.Cells.Item(35, 7).Value  = $productoId
$ws.Cells.Item(35, 8).Value  = $producto
$ws.Cells.Item(35, 9).Value  = $categoriaId
$ws.Cells.Item(35, 10).Value = $categoria
$ws.Cells.Item(35, 11).Value = $variedad
$ws.Cells.Item(35, 12).Value = "Segunda"
$ws.Cells.Item(35, 13).Value = 160
$ws.Cells.Item(35, 14).Value = 14000
$ws.Cells.Item(35, 15).Value = 15000
$ws.Cells.Item(35, 16).Value = 14500
$ws.Cells.Item(35, 17).Value = $unidad
$ws.Cells.Item(35, 18).Value = $origen
$ws.Cells.Item(35, 19).Value = 725
$ws.Cells.Item(35, 20).Value = $kgUnidad
